# Adding synonyms to calcaneus terms (trait sheet), per commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trait")

# Row 2 (calcaneus distal breadth): replace the old narrow "calcaneus width ML"
# synonym with the fuller plantigrade / non-plantigrade synonym list.
$ws.Range("D2").Value = "calcaneus distal width ML (non-plantigrade); calcaneus anteiror width ML (plantigrade); calcaneus anterior breadth (plantigrade)"

# Row 4 (calcaneus proximal length): add a new synonym.
$ws.Range("D4").Value = "calcaneus posterior length (plantigrade)"

# Row 6 (calcaneus proximal breadth): add a new synonym.
$ws.Range("D6").Value = "calcaneus posterior breadth (plantigrade); calcaneus posterior width ML (plantigrade); calcaneus proximal width ML"

# Row 7: TERM renamed from "calcaneus maximal depth" to "calcaneus proximal depth",
# plus a new synonym.
$ws.Range("C7").Value = "calcaneus proximal depth"
$ws.Range("D7").Value = "calcaneus proximal width AP (non-plantigrade); calcaneus posterior depth (plantigrade); calcaneus posterior width PD (plantigrade)"

# Row 8 (calcanus distal depth): replace the old narrow "calcaneus lateral depth"
# synonym with the fuller plantigrade / non-plantigrade synonym list.
$ws.Range("D8").Value = "calcaneus lateral depth; calcaneus distal width AP (non-plantigrade); calcaneus anterior depth (plantigrade); calcneus anterior width PD (plantigrade)"

# Row 9 (depth of calcaneal body): add a new synonym.
$ws.Range("D9").Value = "width AP of calcaneal body (non-plantigrade); width PD of calcaneal body (plantigrade)"

# The TERM column (C) now has longer entries; best-fit the column width like Excel
# does when a user auto-fits the column after editing.
$ws.Columns.Item(3).AutoFit() | Out-Null

# Restore the selection on the non-active "AB" sheet to A2.
$wsAB = $wb.Worksheets.Item("AB")
$wsAB.Range("A2").Select() | Out-Null

# Re-activate "trait" (it stays the visible/selected tab) and move its
# selection to D12.
$ws.Activate() | Out-Null
$ws.Range("D12").Select() | Out-Null
